# small adjustments to DD and DPE for EPICP to reflect dataset file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F8: "PAL" -> "pal"
$ws.Range("F8").Value = "pal"

# Fill in missing input_variables (column F) for rows 97-108,
# mirroring the dataschema_variable (column B) value for each row.
$rows = 97..108
foreach ($r in $rows) {
    $varName = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 6).Value = $varName
}
